$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header label for the transaction amount column (万元 -> 元)
$ws.Range("C1").Value = "交易金额（元）"

# Update investor name values in column B
$ws.Range("B2").Value = "xxxxx1"
$ws.Range("B3").Value = "yyyyy2"
$ws.Range("B4").Value = "zzzzz3"

# Update the active selection to match the saved view state
$ws.Range("D13").Select()
